# Updated symbol list on Mon Dec 12 17:38:38 UTC 2022 with GitHub Actions
#
# Applies the refreshed crypto price snapshot: most rows only get a new
# "Price" (column D) figure, but rows 42/43 (CEJI <-> BKEXToken) swapped
# their list position, so their Coin/Link/Price/Volume cells are updated
# to reflect the new row contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay a *text* string even though it
# looks numeric (e.g. "0.1640"), without leaving a lasting NumberFormat
# change on the cell (match the look of a plain, never-formatted cell).
function Set-TextValue($range, [string]$value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# --- Column D (Price) updates -------------------------------------------
Set-TextValue "D2"  "276.23"
Set-TextValue "D4"  "6.215"
Set-TextValue "D5"  "0.06191"
Set-TextValue "D6"  "3.576"
Set-TextValue "D7"  "1.514"
Set-TextValue "D8"  "6.545"
Set-TextValue "D9"  "0.8230"
Set-TextValue "D10" "0.1640"
Set-TextValue "D11" "0.08218"
Set-TextValue "D12" "0.03424"
Set-TextValue "D14" "0.09127"
Set-TextValue "D15" "3.771"
Set-TextValue "D16" "0.001625"
Set-TextValue "D17" "0.04702"
Set-TextValue "D18" "0.006277"
Set-TextValue "D19" "0.006140"
Set-TextValue "D21" "0.0001501"
Set-TextValue "D22" "3.729"
Set-TextValue "D23" "2.316"
Set-TextValue "D24" "0.01385"
Set-TextValue "D28" "0.0002738"
Set-TextValue "D40" "0.04669"
Set-TextValue "D41" "0.007033"

# --- Rows 42 & 43: CEJI and BKEXToken swap places -------------------------
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1103"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003202"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining column D (Price) updates -----------------------------------
Set-TextValue "D44" "0.01113"
Set-TextValue "D45" "0.00006416"
Set-TextValue "D47" "0.8455"
Set-TextValue "D48" "0.001384"
Set-TextValue "D49" "0.00001901"
